$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -11.04739999999999
$ws.Range("E5").Value = 12.5282
$ws.Range("E9").Value = 12.82630000000001
$ws.Range("E11").Value = 13.18639999999999
$ws.Range("C21").Value = -13.33710000000001
$ws.Range("E21").Value = 12.91299999999999
$ws.Range("C23").Value = -11.9646
$ws.Range("C25").Value = -11.0904
